$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (Pokemon #307), shifting existing rows 9-12 down to 10-13
$ws.Rows("9").Insert()

# Fill in the newly inserted row with data for ability #307 (unused in-game)
# (order matches the shared-string insertion order seen in the target file: E, D, C, B, A)
$ws.Range("E9").Value = "這個特性編號在遊戲中未被使用，參考: https://wiki.52poke.com/wiki/Talk:%E7%89%B9%E6%80%A7%E5%88%97%E8%A1%A8"
$ws.Range("D9").Value = "<No Data>"
$ws.Range("C9").Value = "？？？"
$ws.Range("B9").Value = "<数据暂缺>"
$ws.Range("A9").Value = 307

# Update the defined name range to cover the new row
foreach ($n in $wb.Names) {
    if ($n.Name -eq "工作表1!_1") {
        $n.RefersTo = "=工作表1!`$A`$1:`$G`$13"
    }
}

# Move/restore selection to the newly edited cell
$ws.Range("B9").Select() | Out-Null
